$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 09:32"

# --- Straightforward numeric updates (country stays in the same row) ---

# Row 39: Ucrania
$ws.Range("B39").Value = 54133
$ws.Range("C39").Value = 612
$ws.Range("D39").Value = 26503
$ws.Range("E39").Value = 26232
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 15
$ws.Range("H39").Value = 1398

# Row 53: Armenia
$ws.Range("B53").Value = 32151
$ws.Range("C53").Value = 182
$ws.Range("D53").Value = 19865
$ws.Range("E53").Value = 11713
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 8
$ws.Range("H53").Value = 573

# Row 99: Hungria
$ws.Range("B99").Value = 4247
$ws.Range("C99").Value = 13
$ws.Range("D99").Value = 3073
$ws.Range("E99").Value = 579

# Row 121: Lituania
$ws.Range("B121").Value = 1874
$ws.Range("C121").Value = 5
$ws.Range("E121").Value = 224

# Row 136: Letonia
$ws.Range("B136").Value = 1174
$ws.Range("C136").Value = 1
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 31

# --- Countries whose case counts moved them past a neighbouring row in the
# sorted list. The row data (B:H) is swapped between the two adjacent rows,
# with the new/updated figures landing on whichever row now sorts first. ---

# Rows 143/144: Georgia overtakes Uruguay (Georgia: 995 cases, Uruguay: 987 cases)
$ws.Range("A143").Value = "Georgia"
$ws.Range("B143").Value = 995
$ws.Range("C143").Value = 9
$ws.Range("D143").Value = 857
$ws.Range("E143").Value = 123
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 15

$ws.Range("A144").Value = "Uruguay"
$ws.Range("B144").Value = 987
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 896
$ws.Range("E144").Value = 60
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 31

# Rows 160/161: Siria overtakes Botsuana (Siria: 417 cases, Botsuana: 399 cases)
$ws.Range("A160").Value = "Siria"
$ws.Range("B160").Value = 417
$ws.Range("C160").Value = 23
$ws.Range("D160").Value = 136
$ws.Range("E160").Value = 262
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 3
$ws.Range("H160").Value = 19

$ws.Range("A161").Value = "Botsuana"
$ws.Range("B161").Value = 399
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 38
$ws.Range("E161").Value = 360
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 1

# Rows 209/210: Groenlandia and Islas Malvinas are tied (13 cases each) but
# swap places in the sort order; underlying figures are unchanged.
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
